$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Provider" data rows (10-13), matching the look of the
#     existing "Login" rows (2-5): copy the format of A2 across the
#     whole new block so every new cell (A-F) shares one style, then
#     fill in the text values.
$ws.Range("A2").Copy()
$ws.Range("A10:F13").PasteSpecial(-4122)

$ws.Range("A10").Value = "Provider"
$ws.Range("B10").Value = "POST"

$ws.Range("A11").Value = "Provider"
$ws.Range("B11").Value = "GET"

$ws.Range("A12").Value = "Provider"
$ws.Range("B12").Value = "PUT"

$ws.Range("A13").Value = "Provider"
$ws.Range("B13").Value = "DELETE"

# --- Re-color the header row (row 1): yellow text on a black/dark
#     green fill instead of black text on an orange/amber fill.
$ws.Range("A1:F1").Font.Color = 65535
$ws.Range("A1:F1").Interior.Color = 0
$ws.Range("A1:F1").Interior.PatternColor = 13056

# --- Move the active selection like the source edit did.
$ws.Range("C15").Select() | Out-Null
